$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 165
$ws.Range("A165").Value = 164
$ws.Range("B165").Value = "Looking for the Jester"
$ws.Range("C165").Value = "JesterofTime"
$ws.Range("N165").Value = 1000
$ws.Range("O165").Value = 500
$ws.Range("P165").Value = 10000
$ws.Range("Q165").Value = 250
$ws.Range("T165").Value = 1
$ws.Range("W165").Value = "Delusional Memories"
$ws.Range("X165").Value = 1
$ws.Range("AA165").Value = "There are stories of a jester who speaks of the magic of time. He has the ability to show you the truth of the situation, one that has become foggy in recent days, as if you cannot remember if you visited a realm of ice when others ask.<br /> <br /> You sit in the pub, listening to The Soldier regale you with his tales of the frozen land, where he met his past face-to-face. He tells you of his travels to a realm where the world is a set of memories—memories that are faded, twisted, broken, and damaged.<br /> <br /> Something about his tale sounds off, as if you were the only one you could envision being a part of the stories he tells—as the main character.<br /> <br /> Either way, you listen on, waiting for the hints and stories of a jester. Alas, he continues on his ramblings, although he does drop one tiny piece that almost slips through your bored ears: Nameless Ones."
$ws.Range("AB165").Value = "You have spent some time investigating a new world that has mysteriously appeared as a magical portal opened on the Surface. The Poet and The Wandering Merchant were seen nearby in the last few days, investigating the portal.<br /> <br /> Alas, you managed to sneak through and collect something for The Soldier: a vial of dirt. A single vial of dirt from an old ruin. You never stayed long enough to properly investigate, but something about this new plane of existence seems off."
$ws.Range("AE165").Value = 8

# Row 166
$ws.Range("A166").Value = 165
$ws.Range("B166").Value = "Answers are needed"
$ws.Range("C166").Value = "ThePoet"
$ws.Range("H166").Value = "Federation Broach"
$ws.Range("N166").Value = 1000
$ws.Range("O166").Value = 500
$ws.Range("P166").Value = 10000
$ws.Range("Q166").Value = 250
$ws.Range("U166").Value = "Looking for the Jester"
$ws.Range("W166").Value = "Delusional Memories"
$ws.Range("X166").Value = 2
$ws.Range("AA166").Value = "You find The Poet in the market at Smugglers Port.<br /> <br /> “Child, it has been some time since we last spoke.”<br /> <br /> He seems older and more frail than before. His youth has escaped.<br /> <br /> You tell him about the tales that The Soldier has told you and how they felt off—like they were your tale but weren’t. You wonder if it has to do with the new Plane you recently visited, where the dirt lay that you collected.<br /> <br /> “The tales you hear are from another reality, parallel to our own. It is said in the other reality that the tales you hear in this reality are the ones you have personally lived in another reality, and that eventually the rivers of time all come to a head as they, like all of us, head to the oceans to be absorbed into their ever-flowing vastness.”<br /> <br /> Another reality. Does it have anything to do with the tale of a jester who is involved in time? Or a Nameless one?<br /> <br /> The Poet smiles and moves past you to the next stall to investigate what’s being sold. He turns to you after a moment, in the bustle of the busy market, and states:<br /> <br /> “There is much you know and do not know at the same time. I can tell you about the puppets: the Nameless Ones pull the strings of the puppets. What I can tell you is that a jester, as you call him in your tales, and a woman—a madam—are...”<br /> <br /> He pauses and looks back at an emerald necklace. He picks it up and tells the merchant he will purchase it. After the exchange, he starts to walk on, and you call after him.<br /> <br /> “Are your answers” He calls back, never looking back.<br /> <br /> He disappears into the hustle and bustle of the market. You stand there for a moment, wondering what he meant."
$ws.Range("AB166").Value = "You have been hanging out around Smugglers Port, seeking more answers. Alas, there are none to be had here. You figure it’s time to go to this new plane of existence and figure out who this jester is and what exactly the Nameless Ones are. Are they beyond what you understand The Poet and The Wandering Merchant to be? It has been alluded to in the past that they are cosmic beings. Alas, is there something beyond them?<br /> <br /> It becomes clear to you that your answers will not be found on this plane. You head out, and in the distance, you see the shadow of a wandering man walk by. You think nothing of it. However, it claws at the back of your mind."
$ws.Range("AE166").Value = 8

# Row 167
$ws.Range("A167").Value = 166
$ws.Range("B167").Value = "Why are we naked?"
$ws.Range("C167").Value = "JesterofTime"
$ws.Range("H167").Value = "Clerical Tome of Prairs"
$ws.Range("N167").Value = 1000
$ws.Range("O167").Value = 500
$ws.Range("P167").Value = 10000
$ws.Range("Q167").Value = 250
$ws.Range("U167").Value = "Answers are needed"
$ws.Range("W167").Value = "Delusional Memories"
$ws.Range("X167").Value = 3
$ws.Range("AA167").Value = "“Hello there!” comes a voice.<br /> Y<br /> ou have no idea where you are, but you know you are not on the Surface anymore. You look towards the voice and see a young man dressed in a variety of colors. He is handsome, well-groomed, and his eyes bluer than the sky itself.<br /> <br /> “How are you, child?! Oh, it’s been so long since we had a visitor of these parts. Would you like a tour of our fine establishment?”<br /> <br /> He is full of energy and comes off as extremely charismatic. You explain you are looking for two people: a Nameless One and a jester of some sort.<br /> <br /> “Well, there you go, child! Wrecking all the fun.”<br /> <br /> You look at him with a puzzled look.<br /> <br /> “Does this interest you?”<br /> <br /> He stands fully naked before you. It is quite a sight."
$ws.Range("AB167").Value = "With your back turned to him, you ask him to put clothes on.<br /> <br /> “Why? Am I not attractive enough for you? Don’t you crave me?”<br /> <br /> You tell him he looks no older than sixteen and to dress quickly so we might talk like adults. You also explain that you are not here to indulge in any fantasies, that you are searching for answers, that The Poet sent you.<br /> <br /> “You can turn around. By the way, I am not a juvenile.” His voice seems upset, a tinge of anger.<br /> <br /> You do so, and he is dressed again, still handsome, but dressed. A part of your mind regrets your choices.<br /> <br /> “I can see into the darkest recesses of your mind, child. I am the most beautiful man here. With or without clothes. You will come to see that. As for your precious jester, I know not of what you speak. As for this nameless <br /> nonsense – I know of a woman. A woman in a place where even your memory becomes twisted. That’s it. That’s all I know.”<br /> <br /> He seems so disappointed. He perks up for a moment and states:<br /> <br /> “I do know, I fucking hate The Poet!” He storms off."
$ws.Range("AE167").Value = 8

# Row 168
$ws.Range("A168").Value = 167
$ws.Range("B168").Value = "A cleric in time"
$ws.Range("C168").Value = "WanderingCleric"
$ws.Range("H168").Value = "Broken Childs Doll"
$ws.Range("N168").Value = 1000
$ws.Range("O168").Value = 500
$ws.Range("P168").Value = 10000
$ws.Range("Q168").Value = 250
$ws.Range("U168").Value = "Why are we naked?"
$ws.Range("W168").Value = "Delusional Memories"
$ws.Range("X168").Value = 4
$ws.Range("AA168").Value = "You try to chase after the man who appeared, undressed, and tried to distract you with his beautiful naked body—a sight you do not want to let go of. That’s when you hear that familiar laughter:<br /> <br /> “Ooooooh hoooo hoooo hoooo”<br /> <br /> It echoes in your head, one you have heard before—not from The Queen of Hearts herself in hell, I mean yes, she laughs the same way, but this laughter—this laughter has a sense of darkness about it.<br /> <br /> “Halt, child!”<br /> <br /> A voice stops you in your path, a man holding a mace and dressed in holy garb comes over to you. A knight of The Church? Working for The Federation?<br /> <br /> “What business do you have in these lands, child? Speak not a single lie, for God is watching you and he shall judge your worth based on the weight of your words and the honesty in your heathen heart.”<br /> <br /> What?<br /> <br /> You study the man for a moment, he is also very handsome. Young, like the last man you encountered—who somehow vanished in this desolate wasteland of barren earth.<br /> <br /> You explain you are lost, you encountered a man who was naked and then not naked, who tried to distract you with his “sinfulness” as you put it, almost remembering there is a game to play with these sorts.<br /> <br /> “You found my brother. Where did he go?”<br /> <br /> Your brother?"
$ws.Range("AB168").Value = "You sit by a fire the holy-looking figure has created. He explains he is known as the Wandering Cleric, for he wanders the land searching for his brother. They were separated by their mother. He was always meant to <br /> look after his brother, keep him on the right path, the right direction. Alas, he fell in with a mysterious man, a man who would grant his insatiable desires.<br /> <br /> He continues to explain that his mother is of another plane of existence, trapped there by a man called The Poet, his associate The Wandering Merchant, and another man—a man whose name cannot be uttered.<br /> <br /> His mother would not allow her sons to be separated from her, so in a twist of rage, she separated them from each other, making one of them always long for and search out his lost brother, while the other falls into a pit <br /> of sensual desires and dreams he cannot escape, manipulated by the one thing we all don’t have enough of: Time."
$ws.Range("AE168").Value = 8

# Row 169
$ws.Range("A169").Value = 168
$ws.Range("B169").Value = "The mother and the curse"
$ws.Range("C169").Value = "WanderingCleric"
$ws.Range("H169").Value = "Dwarven Chains"
$ws.Range("N169").Value = 1000
$ws.Range("O169").Value = 500
$ws.Range("P169").Value = 10000
$ws.Range("Q169").Value = 250
$ws.Range("U169").Value = "A cleric in time"
$ws.Range("V169").Value = "Slaves Journal"
$ws.Range("AA169").Value = "You stare at the Wandering Cleric for a moment. This sounds messed up. You explain that to him in such words. He laughs.<br /> <br /> “I always know where he is. He is with: Fliniguss, his lover.”<br /> <br /> And you? You just continue to search for him, your brother. All the while, you know where he is. It confuses you.<br /> <br /> “It is my fate, child. To search out my brother, to the end of time. To search for him. To long for him.”<br /> <br /> In what way, you ask with your brow raised. This story sounds familiar to you. A brother who longs for another. Wasn’t there a soldier and his brother? A sinful, lustful fantasy?<br /> <br /> “Not to that extent, child. I do not suffer my brother's afflictions.”<br /> <br /> You think to yourself, but alas – it rolls right off your tongue: You suffer from another type of delusion, don’t you?<br /> <br /> You’re shocked at your own words. They echo in the silence of the night.<br /> <br /> The Wandering Cleric rises, walks over, and sits beside you. He puts his hand on your lap, a bit far up your lap. He whispers in your ear: “Don’t we all?”<br /> <br /> A shiver runs down your spine. The laughter floods your head: Oooooh hooo hooo hooo."
$ws.Range("AB169").Value = "You awaken in the morning. A naked man lies on his side, full frontal – smiling at you. As you adjust your eyes, you see that it is the charismatic man you met earlier. He sure likes to be naked a lot.<br /> <br /> “Do not worry, child, I am not here for your desires or earthly pleasures. I seek out the love of my life.”<br /> <br /> He stands and dances in the morning light – fully naked, without a care in the world.<br /> <br /> You ask him about Fliniguss, his brother, and his mother. He stops dancing, and as he walks over, it's as if clothing just appears on his body. By the time he sits right in front of you, he is fully clothed.<br /> <br /> “What the heck is wrong with you?” He asks after a moment – irritated at your questions.<br /> <br /> “Can’t you just have some fun? I know you look amazing under all that armor. We don’t have to get nasty, friend, I have someone for that. Oh, the sex – it’s AMAZING!”<br /> <br /> He laughs and falls back to lay on his back.<br /> <br /> “Alas, child, you ask me all the questions I care not for. Whatever will I do with you?”"

# Row 170
$ws.Range("A170").Value = 169
$ws.Range("B170").Value = "The Poet and The Wandering Merchant"
$ws.Range("C170").Value = "ThePoet"
$ws.Range("N170").Value = 1000
$ws.Range("O170").Value = 500
$ws.Range("P170").Value = 10000
$ws.Range("Q170").Value = 250
$ws.Range("U170").Value = "The mother and the curse"
$ws.Range("W170").Value = "Delusional Memories"
$ws.Range("X170").Value = 5
$ws.Range("AA170").Value = "You have spent the last few days searching for the Wandering Cleric and his brother, this man who likes to be naked a lot. He proclaims to be the most beautiful man. One thing does draw your mind back on track from the mysterious and sexy naked man.<br /> <br /> His mother, their mother. Separated – locked away. The Poet, maybe that’s why this mysterious, sexy naked man does not like The Poet. You decide to head back to the Surface and track down the answers you want."
$ws.Range("AB170").Value = "You find him sitting in a pub with The Wandering Merchant. You approach the both of them, and they warmly greet you.<br /> <br /> You sit with them and start to explain your adventures in this new plane.<br /> <br /> “A place of delusional thoughts. A place where people live in a memory they believe to have existed,” states The Wandering Merchant.<br /> <br /> The Poet nods and adds, “A place we thought we locked away all those years ago, back when The Shattering happened.”<br /> <br /> The Shattering sounds familiar, but you cannot place it; it’s almost foggy. The Wandering Merchant sees this confusion on your face and chimes in, almost clearing the fog away.<br /> <br /> “The Shattering happened back when The Child, or The Creator, as you know him, failed to ascend. Failed to ‘move on,’ if you will.” He pauses, and The Poet looks away, refusing to meet either gaze.<br /> <br /> After a moment, The Poet states, “It was her, that I failed to listen to. She knew everything that would happen. I assumed when we closed those doors. Closed her doors, that her own twisted memories would just reside within the prison. Alas, she created her own world. A world of delusions.”<br /> <br /> That sounds off, and you explain what the naked man said. What his brother said. How the naked man hates you.<br /> <br /> “He hates me because I locked his mother, who, because of his afflictions, cursed her children to never be together. To always be alone. She assumed they were together. She could not have that. Her children must be pure. It was only later she discovered the truth, alas, she has been trapped in her own twisted memories.”<br /> <br /> The Poet coughs and tries to catch his breath. The Wandering Merchant places his hand on The Poet’s shoulder as if to ease his friend’s suffering."
$ws.Range("AE170").Value = 8
